$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text value looks like a plain number need to be
# forced to Text format first, otherwise Excel will silently convert them
# to a numeric value (losing formatting such as trailing zeros).
$textCells = @("D5","D13","D15","D16","D20","D21","D22","D25","D27","D28","D30","D34","D35","D37","D42","D43","D46","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.187.75"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.857.83"
$ws.Range("E3").Value = "  +1.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.60%  "

# Row 5 - BNB
$ws.Range("D5").Value = "239.10"
$ws.Range("E5").Value = "  +3.44%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.73%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.58%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  +6.13%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.40%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.29%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.126.86"
$ws.Range("E12").Value = "  +1.56%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "11.49"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14 - was Polygon, now WrappedEther
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.852.17"
$ws.Range("E14").Value = "  +1.82%  "

# Row 15 - was WrappedEther, now Polygon
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.676"
$ws.Range("E15").Value = "  +0.88%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "4.72"
$ws.Range("E16").Value = "  +2.19%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "35.170.53"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18 - Litecoin
$ws.Range("E18").Value = "  +0.44%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +1.21%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "240.78"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +0.55%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +1.50%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.43%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.94%  "

# Row 25 - Monero
$ws.Range("D25").Value = "169.46"
$ws.Range("E25").Value = "  -1.08%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +26.01%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "8.03"
$ws.Range("E27").Value = "  +3.74%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "17.67"
$ws.Range("E28").Value = "  +1.98%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0560"
$ws.Range("E30").Value = "  +1.81%  "

# Row 31 - BinanceUSD
$ws.Range("E31").Value = "  +0.63%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.08%  "

# Row 33 - WEMIXToken
$ws.Range("E33").Value = "  +26.86%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "4.02"
$ws.Range("E34").Value = "  +2.21%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "2.05"
$ws.Range("E35").Value = "  +10.69%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +17.36%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +7.01%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +4.85%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +4.31%  "

# Row 40 - Aave
$ws.Range("E40").Value = "  -1.30%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.348.95"
$ws.Range("E41").Value = "  +0.86%  "

# Row 42 - Kaspa
$ws.Range("D42").Value = "0.0600"
$ws.Range("E42").Value = "  +15.20%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "14.94"
$ws.Range("E43").Value = "  +2.91%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  +2.68%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  +0.34%  "

# Row 46 - Gas
$ws.Range("D46").Value = "12.40"
$ws.Range("E46").Value = "  +42.79%  "

# Row 47 - MXToken
$ws.Range("E47").Value = "  -0.55%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "6.56"
$ws.Range("E48").Value = "  +5.28%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.045.26"
$ws.Range("E49").Value = "  +1.88%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.0683"
$ws.Range("E50").Value = "  +1.80%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value = "  +0.59%  "
